$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text or $null, new Volume(1h) (E) percent text or $null
$updates = @(
    @{ Row = 2; Price = "26.706.40"; Volume = "+0.25%" }
    @{ Row = 3; Price = "1.597.80"; Volume = "+0.34%" }
    @{ Row = 5; Price = "211.44"; Volume = "+0.46%" }
    @{ Row = 6; Price = $null; Volume = "-0.68%" }
    @{ Row = 8; Price = $null; Volume = "+0.58%" }
    @{ Row = 9; Price = $null; Volume = "+1.13%" }
    @{ Row = 10; Price = "19.50"; Volume = "+0.35%" }
    @{ Row = 11; Price = $null; Volume = "+0.13%" }
    @{ Row = 12; Price = "1.822.41"; Volume = "+0.43%" }
    @{ Row = 13; Price = "1.635.12"; Volume = "+2.78%" }
    @{ Row = 14; Price = "4.04"; Volume = "+0.54%" }
    @{ Row = 15; Price = $null; Volume = "+0.46%" }
    @{ Row = 16; Price = "65.32"; Volume = $null }
    @{ Row = 17; Price = "0.0₃0766"; Volume = "+5.35%" }
    @{ Row = 18; Price = "26.671.01"; Volume = "+0.16%" }
    @{ Row = 19; Price = "209.55"; Volume = "+1.25%" }
    @{ Row = 20; Price = $null; Volume = "+0.21%" }
    @{ Row = 21; Price = "7.12"; Volume = "+4.95%" }
    @{ Row = 22; Price = $null; Volume = "+1.27%" }
    @{ Row = 23; Price = $null; Volume = "+0.83%" }
    @{ Row = 24; Price = $null; Volume = "+1.16%" }
    @{ Row = 25; Price = "143.27"; Volume = "-1.65%" }
    @{ Row = 26; Price = "1.01"; Volume = "+0.20%" }
    @{ Row = 27; Price = "7.15"; Volume = "-0.45%" }
    @{ Row = 28; Price = $null; Volume = "+0.22%" }
    @{ Row = 29; Price = "15.33"; Volume = "+0.66%" }
    @{ Row = 30; Price = "0.0519"; Volume = "+3.08%" }
    @{ Row = 31; Price = "1.16"; Volume = "+0.31%" }
    @{ Row = 32; Price = $null; Volume = "+0.82%" }
    @{ Row = 33; Price = $null; Volume = "+1.90%" }
    @{ Row = 34; Price = "1.286.50"; Volume = "+0.53%" }
    @{ Row = 35; Price = $null; Volume = "-6.21%" }
    @{ Row = 36; Price = $null; Volume = "-0.40%" }
    @{ Row = 37; Price = $null; Volume = "+0.32%" }
    @{ Row = 38; Price = $null; Volume = "+0.11%" }
    @{ Row = 39; Price = $null; Volume = "+16.45%" }
    @{ Row = 40; Price = "0.827"; Volume = "-1.15%" }
    @{ Row = 41; Price = "5.45"; Volume = "+0.75%" }
    @{ Row = 42; Price = $null; Volume = "-0.39%" }
    @{ Row = 43; Price = "0.784"; Volume = "-0.26%" }
    @{ Row = 44; Price = "63.36"; Volume = "+0.18%" }
    @{ Row = 45; Price = "1.734.92"; Volume = "+0.40%" }
    @{ Row = 46; Price = "91.36"; Volume = "+1.89%" }
    @{ Row = 47; Price = $null; Volume = "-2.20%" }
    @{ Row = 48; Price = $null; Volume = "+0.35%" }
    @{ Row = 49; Price = "0.0508"; Volume = "+0.46%" }
    @{ Row = 50; Price = $null; Volume = "+0.07%" }
    @{ Row = 51; Price = "7.33"; Volume = "-1.68%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Cells.Item($u.Row, 4)
        if ($u.Price -match "^-?\d+(\.\d+)?$") {
            # Numeric-looking text (single decimal point) must be forced to stay text,
            # matching the original inline-string "Price" column cells.
            $priceCell.NumberFormat = "@"
            $priceCell.Value = $u.Price
            $priceCell.ClearFormats()
        } else {
            $priceCell.Value = $u.Price
        }
    }
    if ($null -ne $u.Volume) {
        $ws.Cells.Item($u.Row, 5).Value = "  " + $u.Volume + "  "
    }
}
